$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AngularMembers")

# Apply an AutoFilter over the data range, filtering column C ("Has Winding")
# down to rows whose value is 1 (mirrors picking "1" in the AutoFilter
# dropdown for that column in the Excel UI).
$rng = $ws.Range("A1:C26")
[void]$rng.AutoFilter(3, @("1"), 7)

# Correct the "Has Winding" flag for IsPointNearLine (row 7) from 1 to 0.
# This happens after the filter was applied, so the row (still matching the
# filter criteria at apply-time) stays visible/un-hidden.
$ws.Range("C7").Value = 0

# Leave the view scrolled/selected where the user ended up after filtering.
[void]$ws.Range("B11").Select()
